# Update simulated team transition-matrix probabilities (Weber St. A)
# Reflects additional simulated games / refreshed Markov-chain estimates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value2 = 0.1612903225806452
$ws.Range("C2").Value2 = 0.6236559139784946
$ws.Range("J2").Value2 = 0.01433691756272401
$ws.Range("P2").Value2 = 0.0931899641577061
$ws.Range("S2").Value2 = 0.1075268817204301

# Row 3
$ws.Range("C3").Value2 = 0.03225806451612903
$ws.Range("J3").Value2 = 0.02150537634408602
$ws.Range("P3").Value2 = 0.7849462365591398
$ws.Range("S3").Value2 = 0.1612903225806452

# Row 4
$ws.Range("J4").Value2 = 0.06818181818181818
$ws.Range("P4").Value2 = 0.6363636363636364
$ws.Range("S4").Value2 = 0.2954545454545455

# Row 6
$ws.Range("B6").Value2 = 0.09004739336492891
$ws.Range("F6").Value2 = 0.05687203791469194
$ws.Range("J6").Value2 = 0.2464454976303317
$ws.Range("O6").Value2 = 0.02369668246445497
$ws.Range("Q6").Value2 = 0.1706161137440758
$ws.Range("R6").Value2 = 0.07582938388625593
$ws.Range("S6").Value2 = 0.3364928909952606

# Row 7
$ws.Range("B7").Value2 = 0.1156462585034014
$ws.Range("D7").Value2 = 0.0272108843537415
$ws.Range("F7").Value2 = 0.06122448979591837
$ws.Range("J7").Value2 = 0.1020408163265306
$ws.Range("Q7").Value2 = 0.1496598639455782
$ws.Range("R7").Value2 = 0.07482993197278912
$ws.Range("S7").Value2 = 0.4693877551020408

# Row 8
$ws.Range("B8").Value2 = 0.09832134292565947
$ws.Range("D8").Value2 = 0.02877697841726619
$ws.Range("F8").Value2 = 0.07434052757793765
$ws.Range("J8").Value2 = 0.09112709832134293
$ws.Range("O8").Value2 = 0.02158273381294964
$ws.Range("Q8").Value2 = 0.1678657074340528
$ws.Range("R8").Value2 = 0.1079136690647482
$ws.Range("S8").Value2 = 0.4100719424460432

# Row 9
$ws.Range("B9").Value2 = 0.09937888198757763
$ws.Range("D9").Value2 = 0.01863354037267081
$ws.Range("F9").Value2 = 0.04968944099378882
$ws.Range("J9").Value2 = 0.1180124223602484
$ws.Range("O9").Value2 = 0.0124223602484472
$ws.Range("Q9").Value2 = 0.2236024844720497
$ws.Range("R9").Value2 = 0.09937888198757763
$ws.Range("S9").Value2 = 0.3788819875776397

# Row 10
$ws.Range("B10").Value2 = 0.1212669683257919
$ws.Range("D10").Value2 = 0.02352941176470588
$ws.Range("F10").Value2 = 0.0751131221719457
$ws.Range("J10").Value2 = 0.08054298642533937
$ws.Range("O10").Value2 = 0.01900452488687783
$ws.Range("Q10").Value2 = 0.2235294117647059
$ws.Range("R10").Value2 = 0.09954751131221719
$ws.Range("S10").Value2 = 0.3574660633484163

# Row 11
$ws.Range("G11").Value2 = 0.1125541125541126
$ws.Range("J11").Value2 = 0.1428571428571428
$ws.Range("K11").Value2 = 0.1991341991341991
$ws.Range("L11").Value2 = 0.5367965367965368
$ws.Range("S11").Value2 = 0.008658008658008658

# Row 12
$ws.Range("G12").Value2 = 0.8095238095238095
$ws.Range("J12").Value2 = 0.1825396825396825
$ws.Range("L12").Value2 = 0
$ws.Range("S12").Value2 = 0.007936507936507936

# Row 13
$ws.Range("G13").Value2 = 0.7142857142857143
$ws.Range("J13").Value2 = 0.2857142857142857
$ws.Range("S13").Value2 = 0

# Row 15
$ws.Range("F15").Value2 = 0.01913875598086124
$ws.Range("H15").Value2 = 0.1722488038277512
$ws.Range("I15").Value2 = 0.07177033492822966
$ws.Range("J15").Value2 = 0.3301435406698565
$ws.Range("K15").Value2 = 0.1100478468899522
$ws.Range("M15").Value2 = 0.004784688995215311
$ws.Range("O15").Value2 = 0.06698564593301436
$ws.Range("S15").Value2 = 0.2248803827751196

# Row 16
$ws.Range("F16").Value2 = 0.005263157894736842
$ws.Range("H16").Value2 = 0.2421052631578947
$ws.Range("I16").Value2 = 0.05263157894736842
$ws.Range("J16").Value2 = 0.4210526315789473
$ws.Range("K16").Value2 = 0.04736842105263158
$ws.Range("M16").Value2 = 0.03157894736842105
$ws.Range("O16").Value2 = 0.06315789473684211
$ws.Range("S16").Value2 = 0.1368421052631579

# Row 17
$ws.Range("F17").Value2 = 0.01466992665036675
$ws.Range("H17").Value2 = 0.2004889975550122
$ws.Range("I17").Value2 = 0.08068459657701711
$ws.Range("J17").Value2 = 0.449877750611247
$ws.Range("K17").Value2 = 0.08312958435207823
$ws.Range("M17").Value2 = 0.019559902200489
$ws.Range("O17").Value2 = 0.05867970660146699
$ws.Range("S17").Value2 = 0.09290953545232274

# Row 18
$ws.Range("F18").Value2 = 0.02010050251256281
$ws.Range("H18").Value2 = 0.1909547738693467
$ws.Range("I18").Value2 = 0.1105527638190955
$ws.Range("J18").Value2 = 0.4623115577889447
$ws.Range("K18").Value2 = 0.07035175879396985
$ws.Range("M18").Value2 = 0.01005025125628141
$ws.Range("O18").Value2 = 0.05527638190954774
$ws.Range("S18").Value2 = 0.08040201005025126

# Row 19
$ws.Range("F19").Value2 = 0.0228789323164919
$ws.Range("H19").Value2 = 0.21163012392755
$ws.Range("I19").Value2 = 0.07816968541468065
$ws.Range("J19").Value2 = 0.3860819828408008
$ws.Range("K19").Value2 = 0.09914204003813155
$ws.Range("M19").Value2 = 0.02383222116301239
$ws.Range("N19").Value2 = 0.0009532888465204957
$ws.Range("O19").Value2 = 0.08579599618684461
$ws.Range("S19").Value2 = 0.09151572926596759
